$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text representation
# instead of being auto-converted to a number by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '35.622.48'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').Value = '1.896.08'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.69%  '
$ws.Range('D5').Value = '247.63'
$ws.Range('E5').Value = '  -3.31%  '
$ws.Range('D6').Value = '0.692'
$ws.Range('E6').Value = '  -5.31%  '
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').Value = '44.14'
$ws.Range('E8').Value = '  +8.58%  '
$ws.Range('E9').Value = '  -4.36%  '
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').Value = '13.16'
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('D13').Value = '2.170.79'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').Value = '0.732'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '1.874.85'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('D17').Value = '35.611.70'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '73.89'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').Value = '0.0₃0825'
$ws.Range('E19').Value = '  -2.89%  '
$ws.Range('D20').Value = '247.32'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('E22').Value = '  -2.88%  '
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('E24').Value = '  +4.55%  '
$ws.Range('E25').Value = '  -9.75%  '
$ws.Range('D26').Value = '166.19'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -2.18%  '
$ws.Range('D28').Value = '18.42'
$ws.Range('E28').Value = '  -2.11%  '
$ws.Range('E29').Value = '  -4.15%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  +7.73%  '
$ws.Range('D32').Value = '4.26'
$ws.Range('E32').Value = '  -2.98%  '
$ws.Range('D33').Value = '0.0582'
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('E35').Value = '  -0.72%  '
$ws.Range('E36').Value = '  -6.21%  '
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('E38').Value = '  -21.91%  '
$ws.Range('D39').Value = '0.0695'
$ws.Range('E39').Value = '  +6.91%  '
$ws.Range('D40').Value = '17.23'
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').Value = '97.77'
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('D42').Value = '0.0215'
$ws.Range('E42').Value = '  -1.82%  '
$ws.Range('D43').Value = '1.10'
$ws.Range('E43').Value = '  -2.59%  '
$ws.Range('D44').Value = '1.297.08'
$ws.Range('E44').Value = '  -2.92%  '
$ws.Range('E45').Value = '  -3.08%  '
$ws.Range('E46').Value = '  +7.62%  '
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').Value = '12.09'
$ws.Range('E49').Value = '  +3.53%  '
$ws.Range('D50').Value = '43.45'
$ws.Range('E50').Value = '  -4.15%  '
$ws.Range('E51').Value = '  -5.48%  '

# Restore the default (unstyled) cell style so no extra formatting
# is introduced, matching the original workbook styling.
$ws.Range("D2:D51").Style = "Normal"

